# Auto-generated Excel COM-interop script
# Appends new sensor-log rows (2026-02-01, ~14:17-14:18) to the Humidity,
# Temperature, Proximity and Camera sheets, mirroring the automated sensor
# logger that produced the rest of this workbook.
#
# Columns A (plain date text, e.g. '2026-02-01') and, on the Humidity sheet,
# E (percentage text, e.g. '76.9%') would otherwise be auto-converted by Excel
# into a date serial / percentage number, so those cells are briefly switched
# to Text format before the value is entered and switched back to the
# workbook's default Normal style afterwards, keeping every new cell a plain,
# unstyled text value exactly like the existing rows.

$wb = $excel.ActiveWorkbook

# --- Humidity sheet: rows 189-202 (Value column holds percentages like '76.9%') ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$dateCell = $wsHumidity.Cells.Item(189, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(189, 2).Value = '14:17:08'
$wsHumidity.Cells.Item(189, 3).Value = '14:00'
$wsHumidity.Cells.Item(189, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(189, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.9%'
$wsHumidity.Cells.Item(189, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(190, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(190, 2).Value = '14:17:09'
$wsHumidity.Cells.Item(190, 3).Value = '14:00'
$wsHumidity.Cells.Item(190, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(190, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.8%'
$wsHumidity.Cells.Item(190, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(191, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(191, 2).Value = '14:17:28'
$wsHumidity.Cells.Item(191, 3).Value = '14:00'
$wsHumidity.Cells.Item(191, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(191, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.9%'
$wsHumidity.Cells.Item(191, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(192, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(192, 2).Value = '14:17:29'
$wsHumidity.Cells.Item(192, 3).Value = '14:00'
$wsHumidity.Cells.Item(192, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(192, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.9%'
$wsHumidity.Cells.Item(192, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(193, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(193, 2).Value = '14:17:31'
$wsHumidity.Cells.Item(193, 3).Value = '14:00'
$wsHumidity.Cells.Item(193, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(193, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.0%'
$wsHumidity.Cells.Item(193, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(194, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(194, 2).Value = '14:17:32'
$wsHumidity.Cells.Item(194, 3).Value = '14:00'
$wsHumidity.Cells.Item(194, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(194, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.9%'
$wsHumidity.Cells.Item(194, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(195, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(195, 2).Value = '14:17:33'
$wsHumidity.Cells.Item(195, 3).Value = '14:00'
$wsHumidity.Cells.Item(195, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(195, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.9%'
$wsHumidity.Cells.Item(195, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(196, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(196, 2).Value = '14:17:38'
$wsHumidity.Cells.Item(196, 3).Value = '14:00'
$wsHumidity.Cells.Item(196, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(196, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.8%'
$wsHumidity.Cells.Item(196, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(197, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(197, 2).Value = '14:17:43'
$wsHumidity.Cells.Item(197, 3).Value = '14:00'
$wsHumidity.Cells.Item(197, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(197, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.8%'
$wsHumidity.Cells.Item(197, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(198, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(198, 2).Value = '14:17:48'
$wsHumidity.Cells.Item(198, 3).Value = '14:00'
$wsHumidity.Cells.Item(198, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(198, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.8%'
$wsHumidity.Cells.Item(198, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(199, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(199, 2).Value = '14:17:53'
$wsHumidity.Cells.Item(199, 3).Value = '14:00'
$wsHumidity.Cells.Item(199, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(199, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.8%'
$wsHumidity.Cells.Item(199, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(200, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(200, 2).Value = '14:17:58'
$wsHumidity.Cells.Item(200, 3).Value = '14:00'
$wsHumidity.Cells.Item(200, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(200, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.7%'
$wsHumidity.Cells.Item(200, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(201, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(201, 2).Value = '14:18:03'
$wsHumidity.Cells.Item(201, 3).Value = '14:00'
$wsHumidity.Cells.Item(201, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(201, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '76.8%'
$wsHumidity.Cells.Item(201, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

$dateCell = $wsHumidity.Cells.Item(202, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsHumidity.Cells.Item(202, 2).Value = '14:18:08'
$wsHumidity.Cells.Item(202, 3).Value = '14:00'
$wsHumidity.Cells.Item(202, 4).Value = 'Bathroom'
$valueCell = $wsHumidity.Cells.Item(202, 5)
$valueCell.NumberFormat = "@"
$valueCell.Value = '77.8%'
$wsHumidity.Cells.Item(202, 6).Value = 'Active'
$dateCell.Style = "Normal"
$valueCell.Style = "Normal"

# --- Temperature sheet: rows 110-122 (Value column holds e.g. '29.4C', already kept as text) ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
$dateCell = $wsTemperature.Cells.Item(110, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(110, 2).Value = '14:17:09'
$wsTemperature.Cells.Item(110, 3).Value = '14:00'
$wsTemperature.Cells.Item(110, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(110, 5).Value = '29.4C'
$wsTemperature.Cells.Item(110, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(111, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(111, 2).Value = '14:17:10'
$wsTemperature.Cells.Item(111, 3).Value = '14:00'
$wsTemperature.Cells.Item(111, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(111, 5).Value = '29.4C'
$wsTemperature.Cells.Item(111, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(112, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(112, 2).Value = '14:17:29'
$wsTemperature.Cells.Item(112, 3).Value = '14:00'
$wsTemperature.Cells.Item(112, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(112, 5).Value = '29.4C'
$wsTemperature.Cells.Item(112, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(113, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(113, 2).Value = '14:17:30'
$wsTemperature.Cells.Item(113, 3).Value = '14:00'
$wsTemperature.Cells.Item(113, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(113, 5).Value = '29.4C'
$wsTemperature.Cells.Item(113, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(114, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(114, 2).Value = '14:17:31'
$wsTemperature.Cells.Item(114, 3).Value = '14:00'
$wsTemperature.Cells.Item(114, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(114, 5).Value = '29.4C'
$wsTemperature.Cells.Item(114, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(115, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(115, 2).Value = '14:17:32'
$wsTemperature.Cells.Item(115, 3).Value = '14:00'
$wsTemperature.Cells.Item(115, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(115, 5).Value = '29.4C'
$wsTemperature.Cells.Item(115, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(116, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(116, 2).Value = '14:17:34'
$wsTemperature.Cells.Item(116, 3).Value = '14:00'
$wsTemperature.Cells.Item(116, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(116, 5).Value = '29.4C'
$wsTemperature.Cells.Item(116, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(117, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(117, 2).Value = '14:17:38'
$wsTemperature.Cells.Item(117, 3).Value = '14:00'
$wsTemperature.Cells.Item(117, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(117, 5).Value = '29.4C'
$wsTemperature.Cells.Item(117, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(118, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(118, 2).Value = '14:17:43'
$wsTemperature.Cells.Item(118, 3).Value = '14:00'
$wsTemperature.Cells.Item(118, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(118, 5).Value = '29.4C'
$wsTemperature.Cells.Item(118, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(119, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(119, 2).Value = '14:17:48'
$wsTemperature.Cells.Item(119, 3).Value = '14:00'
$wsTemperature.Cells.Item(119, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(119, 5).Value = '29.5C'
$wsTemperature.Cells.Item(119, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(120, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(120, 2).Value = '14:17:53'
$wsTemperature.Cells.Item(120, 3).Value = '14:00'
$wsTemperature.Cells.Item(120, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(120, 5).Value = '29.4C'
$wsTemperature.Cells.Item(120, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(121, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(121, 2).Value = '14:17:58'
$wsTemperature.Cells.Item(121, 3).Value = '14:00'
$wsTemperature.Cells.Item(121, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(121, 5).Value = '29.4C'
$wsTemperature.Cells.Item(121, 6).Value = 'Active'
$dateCell.Style = "Normal"

$dateCell = $wsTemperature.Cells.Item(122, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsTemperature.Cells.Item(122, 2).Value = '14:18:03'
$wsTemperature.Cells.Item(122, 3).Value = '14:00'
$wsTemperature.Cells.Item(122, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(122, 5).Value = '29.5C'
$wsTemperature.Cells.Item(122, 6).Value = 'Active'
$dateCell.Style = "Normal"

# --- Proximity sheet: rows 33-34 (door ENTER/EXIT events) ---
$wsProximity = $wb.Worksheets.Item("Proximity")
$dateCell = $wsProximity.Cells.Item(33, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsProximity.Cells.Item(33, 2).Value = '14:17:10'
$wsProximity.Cells.Item(33, 3).Value = '14:00'
$wsProximity.Cells.Item(33, 4).Value = 'Living Room Main Door'
$wsProximity.Cells.Item(33, 5).Value = 'ENTER'
$wsProximity.Cells.Item(33, 6).Value = 'User ENTERED Living Room Main Door'
$dateCell.Style = "Normal"

$dateCell = $wsProximity.Cells.Item(34, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsProximity.Cells.Item(34, 2).Value = '14:17:28'
$wsProximity.Cells.Item(34, 3).Value = '14:00'
$wsProximity.Cells.Item(34, 4).Value = 'Living Room Main Door'
$wsProximity.Cells.Item(34, 5).Value = 'EXIT'
$wsProximity.Cells.Item(34, 6).Value = 'User EXITED Living Room Main Door'
$dateCell.Style = "Normal"

# --- Camera sheet: row 20 (image captured event) ---
$wsCamera = $wb.Worksheets.Item("Camera")
$dateCell = $wsCamera.Cells.Item(20, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-02-01'
$wsCamera.Cells.Item(20, 2).Value = '14:17:27'
$wsCamera.Cells.Item(20, 3).Value = '14:00'
$wsCamera.Cells.Item(20, 4).Value = 'Living Room Main Door'
$wsCamera.Cells.Item(20, 5).Value = 'Image Captured'
$wsCamera.Cells.Item(20, 6).Value = 'Active'
$dateCell.Style = "Normal"

